$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data layout per row (row 2..97): oldA|newB|newC|newD|newE
# newA = oldA + 1 (the whole dataset shifted forward by one day)
$dataText = @"
46069.01041666666|280|1087|1245.682082257708|1367
46069.02083333334|289|1095|1252.115077440379|1384
46069.03125|282|1144|1258.54807262305|1426
46069.04166666666|280|1145|1264.981067805721|1425
46069.05208333334|281|1124|1256.520175197442|1405
46069.0625|285|1073|1248.059282589162|1358
46069.07291666666|277|1117|1239.598389980883|1394
46069.08333333334|275|1118|1231.137497372604|1393
46069.09375|282|1101|1230.27394296652|1383
46069.10416666666|283|1098|1229.410388560437|1381
46069.11458333334|286|0|1228.546834154353|286
46069.125|287|1089|1227.683279748269|1376
46069.13541666666|291|1093|1237.095052489169|1384
46069.14583333334|285|1141|1246.506825230069|1426
46069.15625|0|1144|1255.91859797097|1144
46069.16666666666|292|1140|1265.33037071187|1432
46069.17708333334|347|1164|1297.359506606251|1511
46069.1875|341|1171|1329.388642500632|1512
46069.19791666666|463|1185|1361.417778395013|1648
46069.20833333334|504|1213|1393.446914289394|1717
46069.21875|572|1194|1564.731475314716|1766
46069.22916666666|559|1182|1736.016036340038|1741
46069.23958333334|469|1115|1907.300597365359|1584
46069.25|486|1135|2078.585158390681|1621
46069.26041666666|641|1422|2129.340811181675|2063
46069.27083333334|592|1393|2180.096463972669|1985
46069.28125|589|1391|2230.852116763663|1980
46069.29166666666|0|1399|2281.607769554656|1399
46069.30208333334|694|1572|2274.146271371438|2266
46069.3125|0|0|2266.684773188219|0
46069.32291666666|0|0|2259.223275005|0
46069.33333333334|0|0|2251.761776821781|0
46069.34375|0|0|2201.413644087483|0
46069.35416666666|0|0|2151.065511353185|0
46069.36458333334|0|0|2100.717378618887|0
46069.375|0|0|2050.369245884589|0
46069.38541666666|0|0|2025.918042475339|0
46069.39583333334|0|0|2001.466839066088|0
46069.40625|0|0|1977.015635656838|0
46069.41666666666|0|0|1952.564432247587|0
46069.42708333334|0|0|1907.940986025652|0
46069.4375|0|0|1863.317539803717|0
46069.44791666666|0|0|1818.694093581782|0
46069.45833333334|0|0|1774.070647359847|0
46069.46875|0|0|1759.390222455906|0
46069.47916666666|0|0|1744.709797551964|0
46069.48958333334|0|0|1730.029372648022|0
46069.5|0|0|1715.34894774408|0
46069.51041666666|0|0|1733.347750253616|0
46069.52083333334|0|0|1751.346552763152|0
46069.53125|0|0|1769.345355272688|0
46069.54166666666|0|0|1787.344157782224|0
46069.55208333334|0|0|1805.401177442841|0
46069.5625|0|0|1823.458197103458|0
46069.57291666666|0|0|1841.515216764075|0
46069.58333333334|0|0|1859.572236424692|0
46069.59375|0|0|1926.521960045203|0
46069.60416666666|0|0|1993.471683665714|0
46069.61458333334|0|0|2060.421407286226|0
46069.625|0|0|2127.371130906737|0
46069.63541666666|0|0|2197.58101498187|0
46069.64583333334|0|0|2267.790899057003|0
46069.65625|0|0|2338.000783132137|0
46069.66666666666|0|0|2408.21066720727|0
46069.67708333334|0|0|2437.571517015362|0
46069.6875|0|0|2466.932366823454|0
46069.69791666666|0|0|2496.293216631546|0
46069.70833333334|0|0|2525.654066439638|0
46069.71875|0|0|2540.906959994933|0
46069.72916666666|0|0|2556.159853550227|0
46069.73958333334|0|0|2571.412747105521|0
46069.75|0|0|2586.665640660815|0
46069.76041666666|0|0|2596.882750656745|0
46069.77083333334|0|0|2607.099860652675|0
46069.78125|0|0|2617.316970648605|0
46069.79166666666|0|0|2627.534080644534|0
46069.80208333334|0|0|2614.949473175722|0
46069.8125|0|0|2602.364865706911|0
46069.82291666666|0|0|2589.780258238099|0
46069.83333333334|0|0|2577.195650769287|0
46069.84375|0|0|2487.308369664661|0
46069.85416666666|0|0|2397.421088560035|0
46069.86458333334|0|0|2307.533807455409|0
46069.875|0|0|2217.646526350783|0
46069.88541666666|0|0|2131.591874352294|0
46069.89583333334|0|0|2045.537222353805|0
46069.90625|0|0|1959.482570355315|0
46069.91666666666|0|0|1873.427918356826|0
46069.92708333334|0|0|1751.045764149976|0
46069.9375|0|0|1628.663609943126|0
46069.94791666666|0|0|1506.281455736276|0
46069.95833333334|0|0|1383.899301529426|0
46069.96875|0|0|1328.656755127787|0
46069.97916666666|0|0|1273.414208726148|0
46069.98958333334|0|0|1218.171662324509|0
46070|0|0|1162.92911592287|0
"@

$lines = $dataText -split "`n"
$r = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $oldA = [double]$parts[0]
    $newA = $oldA + 1
    $newB = [double]$parts[1]
    $newC = [double]$parts[2]
    $newD = [double]$parts[3]
    $newE = [double]$parts[4]

    $ws.Cells.Item($r, 1).Value = $newA
    $ws.Cells.Item($r, 2).Value = $newB
    $ws.Cells.Item($r, 3).Value = $newC
    $ws.Cells.Item($r, 4).Value = $newD
    $ws.Cells.Item($r, 5).Value = $newE

    $r = $r + 1
}

Write-Host ("Updated rows 2.." + ($r - 1))
